# Apply the "definitions" dashboard metadata update:
#  - Convert the comma-separated "possible values" lists to pipe-separated
#    lists (so the individual options can contain commas safely).
#  - Clear out the stray "Date" value left in the "Next updated" row.
#  - Move the active selection to D3 (matches the author's last cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value2  = "Indicator| Dashboard| Statistical report"
$ws.Range("D12").Value2 = "broken down by female & male| female only| male only"
$ws.Range("D13").Value2 = "Age| sex| SIMD| deprivation| religion| sexual orientation| income| long term conditions| urban/rural"
$ws.Range("D14").Value2 = "Scotland| NHS Board| Treatment Centre| Intermediate zone| HSCP| Health board| Local Authority| Alcohol & Drug Partnership| GP Practice|  GP Practice Cluster| Hospital| Prison"
$ws.Range("D17").Value2 = "1 yearly| 2 yearly| 3 yearly| 4 yearly| Adhoc| Monthly| Quarterly| Weekly"

# "Next updated" row no longer carries the leftover "Date" note.
$ws.Range("D20").ClearContents()

# Restore the author's last-saved cell selection.
$ws.Range("D3").Select()
